$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 21
$ws.Cells.Item($row, 1).Value = 17
$ws.Cells.Item($row, 2).Value = 210
$ws.Cells.Item($row, 3).Value = "Query Stats"
$ws.Cells.Item($row, 4).Value = "Top Resource-Intensive Queries"
$ws.Cells.Item($row, 5).Value = "http://BrentOzar.com/go/topqueries"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 5), "http://BrentOzar.com/go/topqueries")

$ws.Cells.Item(20, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$ws.Cells.Item($row, 5).Value = "http://BrentOzar.com/go/topqueries"

$ws.Range("A22").Select() | Out-Null
